$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '59.651.00'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +2.66%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.413.73'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  +0.03%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '550.02'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '136.72'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +2.77%  '
$ws.Range("E7").Value = '  -0.01%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.589'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +4.38%  '
$ws.Range("E9").Value = '  +0.01%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '5.68'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +2.28%  '
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("E12").Value = '  +0.04%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '24.65'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.59%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '2.844.89'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +2.27%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '59.647.00'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +2.79%  '
$ws.Range("E16").Value = '  +0.40%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '2.417.48'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +3.26%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '11.25'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +2.53%  '
$ws.Range("E19").Value = '  +0.84%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '329.17'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.28%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.68'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -3.31%  '
$ws.Range("E22").Value = '  +0.19%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '65.64'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +3.43%  '
$ws.Range("E24").Value = '  +3.12%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '8.55'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +3.91%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  +1.42%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0772'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +4.01%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.76'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -0.01%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '170.24'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.49%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '6.12'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -0.30%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '18.60'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +1.26%  '
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("E38").Value = '  +0.04%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '39.30'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +0.31%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '313.44'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +8.88%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  -0.61%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '138.03'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -3.15%  '
$ws.Range("E44").Value = '  +2.19%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.0517'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("E46").Value = '  +3.18%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.577'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +2.25%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0224'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.80%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.392'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +0.90%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '17.48'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.37%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '11.04'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.26%  '
